{"js": "// Update the room-notice schedule line from \"Third Monday\" to \"Second Monday\"\n// (commit: \"update though Aug 2017\").\n//\n// The text \"Room 400:  Third Monday of every month, 5:30 PM to 7:30 PM\" lives\n// in a single paragraph split across a few runs; the word \"Third\" is unique\n// in the document, so a direct search-and-replace on that word is safe and\n// leaves every other run/format untouched.\nconst body = context.document.body;\n\nconst results = body.search(\"Third\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"Third\" in the document body.');\n}\n\n// Replace just the found word in place; the surrounding run formatting\n// (sz/szCs 32) is inherited automatically since we only touch this range.\nresults.items[0].insertText(\"Second\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Update the room-notice schedule line from \"Third Monday\" to \"Second Monday\"\n# (commit: \"update though Aug 2017\").\n#\n# \"Third\" is unique in the document (it only occurs in the\n# \"Room 400:  Third Monday of every month, ...\" line), so locate it with\n# Find and overwrite just that matched range; the run's existing formatting\n# (sz/szCs 32) carries over automatically since only the matched text range\n# is rewritten.\n$d = $word.ActiveDocument\n\n$findRange = $d.Content\n$find = $findRange.Find\n$find.Text = \"Third\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$find.Execute() | Out-Null\n\nif ($find.Found) {\n    $findRange.Text = \"Second\"\n} else {\n    throw \"Could not find 'Third' in the document.\"\n}\n"}
